$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'329.93"
$ws.Range("E2").Value = "'1.21%"
$ws.Range("D3").Value = "'41.38"
$ws.Range("E3").Value = "'4.51%"
$ws.Range("D4").Value = "'5.651"
$ws.Range("E4").Value = "'-1.69%"
$ws.Range("D5").Value = "'0.08218"
$ws.Range("E5").Value = "'2.33%"
$ws.Range("D6").Value = "'8.757"
$ws.Range("E6").Value = "'1.44%"
$ws.Range("D7").Value = "'2.008"
$ws.Range("E7").Value = "'-1.86%"
$ws.Range("D8").Value = "'4.538"
$ws.Range("E8").Value = "'0.87%"
$ws.Range("D9").Value = "'2.971"
$ws.Range("E9").Value = "'1.39%"
$ws.Range("D10").Value = "'0.9256"
$ws.Range("E10").Value = "'0.36%"
$ws.Range("D11").Value = "'0.1275"
$ws.Range("E11").Value = "'1.70%"
$ws.Range("D12").Value = "'0.1961"
$ws.Range("E12").Value = "'0.28%"
$ws.Range("D13").Value = "'0.09366"
$ws.Range("E13").Value = "'1.56%"
$ws.Range("D14").Value = "'0.03839"
$ws.Range("E14").Value = "'7.52%"
$ws.Range("D15").Value = "'0.1061"
$ws.Range("E15").Value = "'1.10%"
$ws.Range("D16").Value = "'0.001308"
$ws.Range("E16").Value = "'0.99%"
$ws.Range("D17").Value = "'0.006169"
$ws.Range("E17").Value = "'0.13%"
$ws.Range("D19").Value = "'3.446"
$ws.Range("E19").Value = "'2.50%"
$ws.Range("E20").Value = "'-0.06%"
$ws.Range("D21").Value = "'8.313"
$ws.Range("E21").Value = "'-5.03%"
$ws.Range("E22").Value = "'1.61%"
$ws.Range("E23").Value = "'-0.16%"
$ws.Range("D24").Value = "'0.04395"
$ws.Range("E24").Value = "'-0.01%"
$ws.Range("E25").Value = "'-0.29%"
$ws.Range("D26").Value = "'0.004316"
$ws.Range("E26").Value = "'-6.35%"
$ws.Range("E27").Value = "'0.90%"
$ws.Range("D39").Value = "'0.02761"
$ws.Range("E39").Value = "'10.75%"
$ws.Range("D40").Value = "'0.05524"
$ws.Range("E40").Value = "'3.97%"
$ws.Range("D41").Value = "'0.007919"
$ws.Range("E41").Value = "'6.19%"
$ws.Range("D42").Value = "'0.1422"
$ws.Range("E42").Value = "'1.24%"
$ws.Range("D43").Value = "'0.008946"
$ws.Range("E43").Value = "'-9.75%"
$ws.Range("D44").Value = "'0.002141"
$ws.Range("E44").Value = "'1.23%"
$ws.Range("D45").Value = "'0.01186"
$ws.Range("E45").Value = "'6.93%"
$ws.Range("D46").Value = "'0.00006989"
$ws.Range("E46").Value = "'4.59%"
$ws.Range("E47").Value = "'0.04%"
$ws.Range("E48").Value = "'5.07%"
$ws.Range("E49").Value = "'0.05%"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("E51").Value = "'0.04%"
